$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Update percentage contribution values (was 3, now 1.3)
$ws.Range("B2").Value = 1.3
$ws.Range("B3").Value = 1.3

# Update Beitragsbemessungsgrenze AV Ost / West values
$ws.Range("B4").Value = 90600
$ws.Range("B5").Value = 89400

# Update Eintragungsdatum value (was 15.12.2023, now 01.01.2024)
$ws.Range("B6").Value = "01.01.2024"

# Update selection to B4:B5 with active cell B4
$ws.Range("B4:B5").Select()
